$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new back-propagation video link first
$ws.Range("B3").Value = "https://www.youtube.com/watch?v=q0pm3BrIUFo"

# Fix up the "Explanation" -> "Explanations" typo in A3
$ws.Range("A3").Value = "Explanations of Back Propagation"

# Add the other two new resource links
$ws.Range("B4").Value = "http://ocw.mit.edu/courses/electrical-engineering-and-computer-science/6-034-artificial-intelligence-fall-2010/readings/MIT6_034F10_netmath.pdf"
$ws.Range("B5").Value = "https://www.youtube.com/watch?v=Ih5Mr93E-2c&hd=1"

# Column A should auto-fit to the new widest text (best-fit width, in
# character units - closest value reachable through the ColumnWidth setter)
$ws.Columns("A").ColumnWidth = 26.5

# Match the final selection from the diff
$ws.Range("B5").Select() | Out-Null
